# Add data for 2022-09-01: the "through" date moves from August 23 to
# August 24, 2022, which adds one more day of carjacking data across
# several neighborhoods/months (current "August" column plus the
# matching weekday in prior years' Augusts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet title / tab name: "Through 2022-08-23" -> "Through 2022-08-24"
$ws.Name = "Through 2022-08-24"

# Column header text for the current (in-progress) month.
$ws.Range("B1").Value = "August 2022 (through August 24)"

# Updated neighborhood/month counts.
$ws.Range("AP2").Value = 5    # Austin, August 2017: 4 -> 5
$ws.Range("BF2").Value = 5    # Austin, August 2015: 4 -> 5
$ws.Range("J4").Value = 10    # North Lawndale, August 2021: 9 -> 10
$ws.Range("B5").Value = 10    # Garfield Park, August 2022 (thru): 9 -> 10
$ws.Range("AH6").Value = 3    # Englewood, August 2018: 2 -> 3
$ws.Range("B7").Value = 7     # Humboldt Park, August 2022 (thru): 8 -> 7
$ws.Range("AH7").Value = 3    # Humboldt Park, August 2018: 2 -> 3
$ws.Range("AP7").Value = 3    # Humboldt Park, August 2017: 2 -> 3
$ws.Range("AP8").Value = 5    # South Shore, August 2017: 4 -> 5
$ws.Range("AX9").Value = 9    # Chatham, August 2016: 8 -> 9
$ws.Range("B12").Value = 7    # West Town, August 2022 (thru): 6 -> 7
$ws.Range("J12").Value = 4    # West Town, August 2021: 3 -> 4
$ws.Range("Z13").Value = 3    # Roseland, August 2019: 2 -> 3
$ws.Range("AP13").Value = 2   # Roseland, August 2017: 1 -> 2
$ws.Range("AX14").Value = 2   # Woodlawn, August 2016: 1 -> 2
$ws.Range("B22").Value = 1    # Loop, August 2022 (thru): (blank) -> 1
$ws.Range("R24").Value = 2    # Lake View, August 2020: 1 -> 2
$ws.Range("J32").Value = 2    # Calumet Heights, August 2021: 1 -> 2
$ws.Range("R34").Value = 2    # Washington Park, August 2020: 1 -> 2
$ws.Range("AP46").Value = 2   # Kenwood, August 2017: 1 -> 2
$ws.Range("R64").Value = 4    # Bucktown, August 2020: 3 -> 4
$ws.Range("AP66").Value = 2   # Chicago Lawn, August 2017: 1 -> 2
$ws.Range("AX66").Value = 3   # Chicago Lawn, August 2016: 2 -> 3
